$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)   # "总计" summary sheet
$q3 = $wb.Worksheets.Item(2)      # existing "2022-Q3" sheet (fund holdings)

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" sheet so its data is preserved as
#    a standalone "2022-Q3" sheet placed right after the original, then
#    turn the original (still in its original slot) into "2022-Q4" and
#    replace its contents with the new quarter's data.
# ---------------------------------------------------------------------
$q3.Copy([System.Reflection.Missing]::Value, $q3)
$q3copy = $wb.Worksheets.Item(3)

$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

$q4 = $q3

# Match the "2022-Q4" sheet's page setup to the one used on "总计"
# (left/right = 0.75in, top/bottom = 1in, header/footer = 0.5in;
# COM PageSetup margins are expressed in points, 72pt = 1in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Q4 only has two funds on record (vs four for Q3) - drop the extra rows.
$q4.Rows("4:5").Delete()

# Re-use the "总计" header/index formatting (bold + bordered, centred)
# for the "2022-Q4" header row and the numeric index column.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

# Columns that hold numeric-looking text (fund codes, percentages, etc.)
# must be pre-formatted as Text so Excel doesn't coerce them to numbers.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:D3").NumberFormat = "@"
$q4.Range("E2:E3").NumberFormat = "@"
$q4.Range("F2:F3").NumberFormat = "@"
$q4.Range("G2").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "001219"
$q4.Range("C2").Value = "上投摩根动态多因子策略混合A"
$q4.Range("D2").Value = "0.97"
$q4.Range("E2").Value = "92.08"
$q4.Range("F2").Value = "3.78"
$q4.Range("G2").Value = "0.0367"
$q4.Range("H2").Value = 7

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "017176"
$q4.Range("C3").Value = "上投摩根动态多因子策略混合C"
$q4.Range("D3").Value = "0.00"
$q4.Range("E3").Value = "92.08"
$q4.Range("F3").Value = "3.78"
$q4.Range("G3").Value = 0
$q4.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: push the existing 2022-Q3 summary
#    row down to row 3 and put the new 2022-Q4 summary in row 2.
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("B3").Value = $total.Range("B2").Value()
$total.Range("C3").Value = $total.Range("C2").Value()
$total.Range("D3").Value = $total.Range("D2").Value()
$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

Write-Output "done"
